# Generate Report for Handoff
#
# The localization status report is regenerated: the files that were
# previously reported as "Handed back: in sync with en-US" are now
# "Ready for handoff" again, the handoff timestamps are refreshed, and the
# (now shorter) status/datetime columns are narrowed to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ---
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
$wsOverview.Range("G2").Value = "2016-09-03 15:04:21"
$wsDeDe.Range("H2").Value     = "2016-09-03 15:04:21"

# zh-cn "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-03 15:04:14"

# --- Narrow the now-shorter status/datetime columns ---
# (target character width ~17.22; ColumnWidth snaps to the nearest pixel
# boundary, so feed a pre-snapped value to land as close as possible)
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth     = 16.33
$wsDeDe.Range("C1").ColumnWidth     = 16.33
